$d = $word.ActiveDocument

# --- First paragraph: "**ID__AFFARS_5349_topic_5__ID** " -----------------
# The paragraph currently holds two runs: the placeholder-id text, and a
# trailing run that is a single literal space. The edit drops that second
# run entirely and renames the placeholder id.

# 1. Remove the stand-alone trailing-space run (paragraph 1 spans
#    characters 0-32; the id text occupies 0-31, the space run is 31-32).
$spaceRun = $d.Range(31, 32)
$spaceRun.Delete()

# 2. Rename the placeholder id text in what is now the sole run.
$idRun = $d.Range(0, 31)
$idRun.Text = "**ID__AFFARS_SUBPART_5349_4__ID**"

# 3. Paragraph formatting tweaks on paragraph 1.
$p1 = $d.Paragraphs(1)

# Left indent: 120 twips -> 225 twips (Word's LeftIndent is in points).
$p1.Format.LeftIndent = 225 / 20

# Add a thin paragraph border on all four sides with 5pt text spacing.
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
